$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete old column A (the numeric index column) so that B:F shift left to A:E
$ws.Columns.Item(1).Delete()
